# Reorder the "Requisitos" bullet list:
#   before: LOM3011, LOM3036, LOM3046, LOM3013
#   after:  LOM3036, LOM3011, LOM3013, LOM3046
#
# Each requirement line lives in its own <w:r> run (text + line break), all
# with identical/default formatting. We replace the whole paragraph's XML
# in one shot via Range.InsertXML so the run-per-line structure (and the
# w:br at the end of each run) is preserved exactly, just reordered.

$d = $word.ActiveDocument

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "LOM3011*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the Requisitos paragraph (starting with 'LOM3011')."
}

$w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

$xml = '<w:p xmlns:w="' + $w + '">' + `
    '<w:pPr><w:pStyle w:val="ListBullet"/></w:pPr>' + `
    '<w:r><w:t>LOM3036 -  Propriedades Mecânicas  (Requisito)</w:t><w:br/></w:r>' + `
    '<w:r><w:t>LOM3011 -  Ensaios Mecânicos  (Requisito)</w:t><w:br/></w:r>' + `
    '<w:r><w:t>LOM3013 -  Ciência dos Materiais  (Requisito)</w:t><w:br/></w:r>' + `
    '<w:r><w:t>LOM3046 -  Técnicas de Análise Microestrutural  (Requisito)</w:t><w:br/></w:r>' + `
    '</w:p>'

$target.Range.InsertXML($xml)
